$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Small "mistake" corrections on rows 2-4 (Function coefficients) ---
$ws.Range("C2").Value = 16
$ws.Range("C3").Value = 12
$ws.Range("C4").Value = 8

# --- Rename existing shop rows 8-9 and bump their inventory columns (C:F) ---
$ws.Range("A8").Value = "shop0"
$ws.Range("C8").Value = 100000
$ws.Range("D8").Value = 100000
$ws.Range("E8").Value = 100000
$ws.Range("F8").Value = 100000

$ws.Range("A9").Value = "shop1"
$ws.Range("C9").Value = 100000
$ws.Range("D9").Value = 100000
$ws.Range("E9").Value = 100000
$ws.Range("F9").Value = 100000

# --- Add 18 more shop rows (10-27): shop2 .. shop19 ---
for ($i = 2; $i -le 19; $i++) {
    $r = 8 + $i
    $ws.Cells.Item($r, 1).Value = "shop$i"
    $ws.Cells.Item($r, 2).Value = 100
    $ws.Cells.Item($r, 3).Value = 100000
    $ws.Cells.Item($r, 4).Value = 100000
    $ws.Cells.Item($r, 5).Value = 100000
    $ws.Cells.Item($r, 6).Value = 100000

    if ($r -le 17) {
        $ws.Cells.Item($r, 7).Value = 5
    } elseif ($r -le 22) {
        $ws.Cells.Item($r, 7).Value = 7.5
        $ws.Cells.Item($r, 7).NumberFormat = "0.00"
    } else {
        $ws.Cells.Item($r, 7).Value = 10
    }

    $ws.Cells.Item($r, 8).Value = 10
    $ws.Cells.Item($r, 9).Value = 130
    $ws.Cells.Item($r, 10).Value = 6
}

# --- Selection moves to D5 ---
$ws.Range("D5").Select()
